$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D2').Value = '28.689.35'
$ws.Range('E2').Value = '  -3.12%  '
$ws.Range('D3').Value = '1.952.90'
$ws.Range('E3').Value = '  -2.47%  '
$ws.Range('D4').NumberFormat = '@'
$ws.Range('D4').Value = '1.014'
$ws.Range('E4').Value = '  -0.11%  '
$ws.Range('D5').NumberFormat = '@'
$ws.Range('D5').Value = '321.37'
$ws.Range('E6').Value = '  -0.03%  '
$ws.Range('D7').NumberFormat = '@'
$ws.Range('D7').Value = '0.4765'
$ws.Range('E7').Value = '  -4.82%  '
$ws.Range('D8').NumberFormat = '@'
$ws.Range('D8').Value = '0.4025'
$ws.Range('E8').Value = '  -4.85%  '
$ws.Range('D9').NumberFormat = '@'
$ws.Range('D9').Value = '53.72'
$ws.Range('E9').Value = '  -0.40%  '
$ws.Range('D10').NumberFormat = '@'
$ws.Range('D10').Value = '0.08494'
$ws.Range('E10').Value = '  -5.47%  '
$ws.Range('D11').NumberFormat = '@'
$ws.Range('D11').Value = '1.060'
$ws.Range('E11').Value = '  -5.26%  '
$ws.Range('D12').NumberFormat = '@'
$ws.Range('D12').Value = '22.15'
$ws.Range('E12').Value = '  -5.15%  '
$ws.Range('D13').Value = '1.957.64'
$ws.Range('E13').Value = '  -2.65%  '
$ws.Range('D14').NumberFormat = '@'
$ws.Range('D14').Value = '7.604'
$ws.Range('E14').Value = '  -5.60%  '
$ws.Range('D15').NumberFormat = '@'
$ws.Range('D15').Value = '6.207'
$ws.Range('E15').Value = '  -4.19%  '
$ws.Range('D16').NumberFormat = '@'
$ws.Range('D16').Value = '1.015'
$ws.Range('E16').Value = '  -0.02%  '
$ws.Range('D17').NumberFormat = '@'
$ws.Range('D17').Value = '0.00001074'
$ws.Range('E17').Value = '  -3.53%  '
$ws.Range('D18').NumberFormat = '@'
$ws.Range('D18').Value = '89.20'
$ws.Range('E18').Value = '  -5.10%  '
$ws.Range('D19').NumberFormat = '@'
$ws.Range('D19').Value = '0.06614'
$ws.Range('E19').Value = '  -0.75%  '
$ws.Range('D20').NumberFormat = '@'
$ws.Range('D20').Value = '18.64'
$ws.Range('E20').Value = '  -5.59%  '
$ws.Range('E21').Value = '  +0.01%  '
$ws.Range('D22').NumberFormat = '@'
$ws.Range('D22').Value = '5.799'
$ws.Range('E22').Value = '  -2.66%  '
$ws.Range('D23').Value = '28.704.31'
$ws.Range('E23').Value = '  -3.17%  '
$ws.Range('D24').NumberFormat = '@'
$ws.Range('D24').Value = '11.51'
$ws.Range('E24').Value = '  -3.92%  '
$ws.Range('E25').Value = '  +0.10%  '
$ws.Range('D26').Value = '2.199.62'
$ws.Range('E26').Value = '  -2.13%  '
$ws.Range('D27').NumberFormat = '@'
$ws.Range('D27').Value = '154.74'
$ws.Range('E27').Value = '  -3.02%  '
$ws.Range('D28').NumberFormat = '@'
$ws.Range('D28').Value = '20.17'
$ws.Range('E28').Value = '  -2.64%  '
$ws.Range('D29').NumberFormat = '@'
$ws.Range('D29').Value = '5.944'
$ws.Range('E29').Value = '  -7.47%  '
$ws.Range('D30').NumberFormat = '@'
$ws.Range('D30').Value = '2.149'
$ws.Range('E30').Value = '  -6.53%  '
$ws.Range('D31').NumberFormat = '@'
$ws.Range('D31').Value = '123.75'
$ws.Range('E31').Value = '  -3.70%  '
$ws.Range('D32').NumberFormat = '@'
$ws.Range('D32').Value = '1.001'
$ws.Range('E32').Value = '  -4.78%  '
$ws.Range('D33').NumberFormat = '@'
$ws.Range('D33').Value = '0.09569'
$ws.Range('E33').Value = '  -3.72%  '
$ws.Range('D34').NumberFormat = '@'
$ws.Range('D34').Value = '5.653'
$ws.Range('E34').Value = '  -3.05%  '
$ws.Range('E35').Value = '  -3.61%  '
$ws.Range('D36').NumberFormat = '@'
$ws.Range('D36').Value = '1.434'
$ws.Range('E36').Value = '  -8.37%  '
$ws.Range('E37').Value = '  -5.01%  '
$ws.Range('B38').Value = 'Hedera'
$ws.Range('C38').Value = 'https://coinranking.com/coin/jad286TjB+hedera-hbar'
$ws.Range('D38').NumberFormat = '@'
$ws.Range('D38').Value = '0.06214'
$ws.Range('E38').Value = '  -2.20%  '
$ws.Range('B39').Value = 'TrustWalletToken'
$ws.Range('C39').Value = 'https://coinranking.com/coin/Hm3OlynlC+trustwallettoken-twt'
$ws.Range('D39').NumberFormat = '@'
$ws.Range('D39').Value = '1.264'
$ws.Range('E39').Value = '  -3.32%  '
$ws.Range('D40').NumberFormat = '@'
$ws.Range('D40').Value = '8.744'
$ws.Range('E40').Value = '  -6.92%  '
$ws.Range('D41').NumberFormat = '@'
$ws.Range('D41').Value = '0.6208'
$ws.Range('E41').Value = '  -5.52%  '
$ws.Range('D42').NumberFormat = '@'
$ws.Range('D42').Value = '11.07'
$ws.Range('E42').Value = '  -5.28%  '
$ws.Range('D43').NumberFormat = '@'
$ws.Range('D43').Value = '1.013'
$ws.Range('E43').Value = '  -0.02%  '
$ws.Range('D44').NumberFormat = '@'
$ws.Range('D44').Value = '0.1915'
$ws.Range('E44').Value = '  -6.66%  '
$ws.Range('E45').Value = '  +2.42%  '
$ws.Range('D46').NumberFormat = '@'
$ws.Range('D46').Value = '0.5922'
$ws.Range('E46').Value = '  -6.50%  '
$ws.Range('D47').NumberFormat = '@'
$ws.Range('D47').Value = '12.85'
$ws.Range('E47').Value = '  -4.48%  '
$ws.Range('D48').NumberFormat = '@'
$ws.Range('D48').Value = '2.076'
$ws.Range('E48').Value = '  -5.60%  '
$ws.Range('B49').Value = 'PancakeSwap'
$ws.Range('C49').Value = 'https://coinranking.com/coin/ncYFcP709+pancakeswap-cake'
$ws.Range('D49').NumberFormat = '@'
$ws.Range('D49').Value = '3.414'
$ws.Range('E49').Value = '  -3.11%  '
$ws.Range('B50').Value = 'BabyDogeCoin'
$ws.Range('C50').Value = 'https://coinranking.com/coin/JY1_q2c0g+babydogecoin-babydoge'
$ws.Range('D50').NumberFormat = '@'
$ws.Range('D50').Value = '0.00000000334'
$ws.Range('E50').Value = '  -1.32%  '
$ws.Range('D51').NumberFormat = '@'
$ws.Range('D51').Value = '0.06818'
$ws.Range('E51').Value = '  -2.45%  '
